$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.371.44"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "3.513.55"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.21"
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.62"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +6.92%  "
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.389"
$ws.Range("E11").Value = "  +3.83%  "
$ws.Range("D12").Value = "4.112.75"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").Value = "3.513.96"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.85"
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("D17").Value = "64.350.66"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("E19").Value = "  +3.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.57"
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "394.08"
$ws.Range("E21").Value = "  +2.70%  "
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("D23").Value = "3.654.31"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.66"
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000117"
$ws.Range("E27").Value = "  +2.75%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("E32").Value = "  -6.76%  "
$ws.Range("E33").Value = "  +7.23%  "
$ws.Range("D34").Value = "3.544.58"
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.41"
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167.02"
$ws.Range("E40").Value = "  +1.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0790"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.66"
$ws.Range("E43").Value = "  -2.83%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.66"
$ws.Range("E46").Value = "  +2.47%  "
$ws.Range("E47").Value = "  -1.79%  "
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("D49").Value = "2.394.77"
$ws.Range("E49").Value = "  -3.46%  "
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("E51").Value = "  +0.32%  "
